$d = $word.ActiveDocument

function Get-ParaByExactText($doc, [string]$target) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $t = $t.TrimEnd([char]13, [char]7)
        }
        if ($t -eq $target) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) Title paragraph: split "C# Fundamentals for Absolute beginners"
#    into runs separated by spell-check proofErr markers.
# ---------------------------------------------------------------
$pTitle = Get-ParaByExactText $d "C# Fundamentals for Absolute beginners"
if ($pTitle -eq $null) { throw "title paragraph not found" }
$titleXml = @"
<w:p w14:paraId="6F1D432E" w14:textId="77777777" w:rsidR="0001531F" w:rsidRDefault="009A7CAA" w:rsidP="009A7CAA"><w:pPr><w:pStyle w:val="Overskrift1"/></w:pPr><w:r><w:t xml:space="preserve">C# </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fundamentals</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Absolute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>beginners</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
$pTitle.Range.InsertXML($titleXml)

# ---------------------------------------------------------------
# 2) "Use camelCasing for increased readability"
# ---------------------------------------------------------------
$pCamel = Get-ParaByExactText $d "Use camelCasing for increased readability"
if ($pCamel -eq $null) { throw "camelCasing paragraph not found" }
$camelXml = @"
<w:p w14:paraId="6BE14040" w14:textId="76F125D5" w:rsidR="009E0731" w:rsidRDefault="00E01EC2" w:rsidP="009E0731"><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>camelCasing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> for increased readability</w:t></w:r></w:p>
"@
$pCamel.Range.InsertXML($camelXml)

# ---------------------------------------------------------------
# 3) "C# is CaseSensitive!"
# ---------------------------------------------------------------
$pCase = Get-ParaByExactText $d "C# is CaseSensitive!"
if ($pCase -eq $null) { throw "CaseSensitive paragraph not found" }
$caseXml = @"
<w:p w14:paraId="4AE97F6C" w14:textId="18B8190B" w:rsidR="00AE2D32" w:rsidRDefault="008124E4" w:rsidP="00AE2D32"><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">C# is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>CaseSensitive</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>!</w:t></w:r></w:p>
"@
$pCase.Range.InsertXML($caseXml)

# ---------------------------------------------------------------
# 4) "Intellisense may automatically pop up" - wrap "Intellisense" with proofErr
# ---------------------------------------------------------------
$pIntelli = Get-ParaByExactText $d "Intellisense may automatically pop up"
if ($pIntelli -eq $null) { throw "Intellisense paragraph not found" }
$intelliXml = @"
<w:p w14:paraId="558A35C5" w14:textId="77777777" w:rsidR="00586F4D" w:rsidRDefault="00586D87" w:rsidP="00AE2D32"><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Intellisense</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00586F4D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> may automatically pop up</w:t></w:r></w:p>
"@
$pIntelli.Range.InsertXML($intelliXml)

# ---------------------------------------------------------------
# 5) Last paragraph ("One-line-code-block does not need curly braces."
#    plus the _GoBack bookmark) -> strip the bookmark from it and
#    append 5 new list paragraphs, moving the bookmark into the new,
#    now-last, empty paragraph.
# ---------------------------------------------------------------
$pLast = Get-ParaByExactText $d "One-line-code-block does not need curly braces."
if ($pLast -eq $null) { throw "last paragraph not found" }

$lastBlockXml = @"
<w:p w14:paraId="34FBAE9B" w14:textId="3B98121C" w:rsidR="008421BE" w:rsidRPr="00A1520A" w:rsidRDefault="00BD7F3A" w:rsidP="00A1520A"><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>One-line-code-block does not need curly braces.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> += something; is short for: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + something;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>If(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> == value) ? result-if-yes; result-if-false;</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Writeline</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">&#8220;You won a {0}.&#8221;, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>var</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>);  //replace the {0} with the value of var.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="2B91AF"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>Console</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>.WriteLine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>&quot;You entered: {0}, therefor you won a {1}.&quot;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>userValue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t>, message);</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="000000"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> //works with several replacements like this</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
"@
$pLast.Range.InsertXML($lastBlockXml)

Write-Output "All edits applied"
